{"js": "const body = context.document.body;\nconst replacements = [\n  [\"722\u00f77=103, 1\", \"813\u00f79=90, 3\"],\n  [\"878\u00f78=109, 6\", \"258\u00f76=43, 0\"],\n  [\"748\u00f79=83, 1\", \"406\u00f78=50, 6\"],\n  [\"133\u00f72=66, 1\", \"540\u00f75=108, 0\"],\n  [\"167\u00f79=18, 5\", \"973\u00f79=108, 1\"],\n  [\"781\u00f72=390, 1\", \"389\u00f77=55, 4\"],\n  [\"966\u00f72=483, 0\", \"509\u00f76=84, 5\"],\n  [\"419\u00f76=69, 5\", \"735\u00f72=367, 1\"],\n  [\"758\u00f72=379, 0\", \"440\u00f77=62, 6\"],\n  [\"819\u00f79=91, 0\", \"204\u00f78=25, 4\"],\n  [\"807\u00f75=161, 2\", \"101\u00f78=12, 5\"],\n  [\"870\u00f73=290, 0\", \"858\u00f73=286, 0\"],\n  [\"370\u00f73=123, 1\", \"839\u00f78=104, 7\"],\n  [\"229\u00f76=38, 1\", \"819\u00f72=409, 1\"],\n  [\"682\u00f77=97, 3\", \"185\u00f77=26, 3\"],\n  [\"992\u00f77=141, 5\", \"670\u00f74=167, 2\"],\n  [\"108\u00f76=18, 0\", \"282\u00f75=56, 2\"],\n  [\"852\u00f78=106, 4\", \"637\u00f75=127, 2\"],\n  [\"724\u00f78=90, 4\", \"694\u00f77=99, 1\"],\n  [\"706\u00f79=78, 4\", \"784\u00f79=87, 1\"],\n  [\"144\u00f72=72, 0\", \"924\u00f72=462, 0\"],\n  [\"396\u00f72=198, 0\", \"159\u00f77=22, 5\"],\n  [\"284\u00f77=40, 4\", \"206\u00f72=103, 0\"],\n  [\"997\u00f73=332, 1\", \"730\u00f74=182, 2\"],\n  [\"656\u00f73=218, 2\", \"163\u00f72=81, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"722\u00f77=103, 1\", \"813\u00f79=90, 3\")\n    ,@(\"878\u00f78=109, 6\", \"258\u00f76=43, 0\")\n    ,@(\"748\u00f79=83, 1\", \"406\u00f78=50, 6\")\n    ,@(\"133\u00f72=66, 1\", \"540\u00f75=108, 0\")\n    ,@(\"167\u00f79=18, 5\", \"973\u00f79=108, 1\")\n    ,@(\"781\u00f72=390, 1\", \"389\u00f77=55, 4\")\n    ,@(\"966\u00f72=483, 0\", \"509\u00f76=84, 5\")\n    ,@(\"419\u00f76=69, 5\", \"735\u00f72=367, 1\")\n    ,@(\"758\u00f72=379, 0\", \"440\u00f77=62, 6\")\n    ,@(\"819\u00f79=91, 0\", \"204\u00f78=25, 4\")\n    ,@(\"807\u00f75=161, 2\", \"101\u00f78=12, 5\")\n    ,@(\"870\u00f73=290, 0\", \"858\u00f73=286, 0\")\n    ,@(\"370\u00f73=123, 1\", \"839\u00f78=104, 7\")\n    ,@(\"229\u00f76=38, 1\", \"819\u00f72=409, 1\")\n    ,@(\"682\u00f77=97, 3\", \"185\u00f77=26, 3\")\n    ,@(\"992\u00f77=141, 5\", \"670\u00f74=167, 2\")\n    ,@(\"108\u00f76=18, 0\", \"282\u00f75=56, 2\")\n    ,@(\"852\u00f78=106, 4\", \"637\u00f75=127, 2\")\n    ,@(\"724\u00f78=90, 4\", \"694\u00f77=99, 1\")\n    ,@(\"706\u00f79=78, 4\", \"784\u00f79=87, 1\")\n    ,@(\"144\u00f72=72, 0\", \"924\u00f72=462, 0\")\n    ,@(\"396\u00f72=198, 0\", \"159\u00f77=22, 5\")\n    ,@(\"284\u00f77=40, 4\", \"206\u00f72=103, 0\")\n    ,@(\"997\u00f73=332, 1\", \"730\u00f74=182, 2\")\n    ,@(\"656\u00f73=218, 2\", \"163\u00f72=81, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
